$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.041804464500799
$ws.Range("D2").Value = 1.051057722846996
$ws.Range("E2").Value = 1.040025812447301
$ws.Range("F2").Value = 1.059670161041093
$ws.Range("I2").Value = 1.046515349916376
$ws.Range("J2").Value = 1.046883698422676
$ws.Range("K2").Value = 1.053810279920723
$ws.Range("L2").Value = 1.0428093309011
$ws.Range("M2").Value = 1.062399032429219

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.042699779461084
$ws.Range("D3").Value = 1.051785684973364
$ws.Range("E3").Value = 1.040786244702254
$ws.Range("F3").Value = 1.060529093033269
$ws.Range("I3").Value = 1.046786279295741
$ws.Range("J3").Value = 1.047425461809878
$ws.Range("K3").Value = 1.05435103014468
$ws.Range("L3").Value = 1.043380233574392
$ws.Range("M3").Value = 1.06307212643042

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.04327963670386
$ws.Range("D4").Value = 1.052257150708121
$ws.Range("E4").Value = 1.041279114036585
$ws.Range("F4").Value = 1.061085641440234
$ws.Range("I4").Value = 1.046960567263925
$ws.Range("J4").Value = 1.047775896780014
$ws.Range("K4").Value = 1.054700686201352
$ws.Range("L4").Value = 1.043749801297007
$ws.Range("M4").Value = 1.063507779662342

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.043523533902915
$ws.Range("D5").Value = 1.052455455200362
$ws.Range("E5").Value = 1.041486510733213
$ws.Range("F5").Value = 1.061319795130716
$ws.Range("I5").Value = 1.047033592648052
$ws.Range("J5").Value = 1.047923189257279
$ws.Range("K5").Value = 1.0548476216312
$ws.Range("L5").Value = 1.043905203574987
$ws.Range("M5").Value = 1.063690954743023

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.0435644926198
$ws.Range("D6").Value = 1.052488757232648
$ws.Range("E6").Value = 1.041521344906079
$ws.Range("F6").Value = 1.061359121115607
$ws.Range("I6").Value = 1.047045839526733
$ws.Range("J6").Value = 1.047947918497557
$ws.Range("K6").Value = 1.054872289184393
$ws.Range("L6").Value = 1.043931298377178
$ws.Range("M6").Value = 1.063721712155978

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.043282895179825
$ws.Range("D7").Value = 1.052259800070339
$ws.Range("E7").Value = 1.041281884519432
$ws.Range("F7").Value = 1.061088769505224
$ws.Range("I7").Value = 1.046961543997376
$ws.Range("J7").Value = 1.047777865030117
$ws.Range("K7").Value = 1.054702649796931
$ws.Range("L7").Value = 1.043751877649903
$ws.Range("M7").Value = 1.06351022715452

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.042106930131979
$ws.Range("D8").Value = 1.051303651989023
$ws.Range("E8").Value = 1.040282633586254
$ws.Range("F8").Value = 1.059960282349027
$ws.Range("I8").Value = 1.046607122823236
$ws.Range("J8").Value = 1.047066814507837
$ws.Range("K8").Value = 1.053993079098522
$ws.Range("L8").Value = 1.043002237195611
$ws.Range("M8").Value = 1.062626482750753

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.040038836832394
$ws.Range("D9").Value = 1.049622126129069
$ws.Range("E9").Value = 1.038528167532488
$ws.Range("F9").Value = 1.057977653747418
$ws.Range("I9").Value = 1.045974796983972
$ws.Range("J9").Value = 1.045812967223023
$ws.Range("K9").Value = 1.052740897854565
$ws.Range("L9").Value = 1.041682520098875
$ws.Range("M9").Value = 1.061070165836894

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.038662948291844
$ws.Range("D10").Value = 1.048503441270084
$ws.Range("E10").Value = 1.037362882314298
$ws.Range("F10").Value = 1.056659974472931
$ws.Range("I10").Value = 1.045548052396475
$ws.Range("J10").Value = 1.044976541133953
$ws.Range("K10").Value = 1.051904953384449
$ws.Range("L10").Value = 1.040803618527564
$ws.Range("M10").Value = 1.060033345802535

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.038067863425091
$ws.Range("D11").Value = 1.048019611322684
$ws.Range("E11").Value = 1.036859354773077
$ws.Range("F11").Value = 1.056090391787377
$ws.Range("I11").Value = 1.045362044470578
$ws.Range("J11").Value = 1.044614248293595
$ws.Range("K11").Value = 1.051542720730963
$ws.Range("L11").Value = 1.040423275478066
$ws.Range("M11").Value = 1.059584580615279

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.037846926454083
$ws.Range("D12").Value = 1.047839982251594
$ws.Range("E12").Value = 1.036672481379919
$ws.Range("F12").Value = 1.055878972495927
$ws.Range("I12").Value = 1.045292769557614
$ws.Range("J12").Value = 1.044479660468145
$ws.Range("K12").Value = 1.051408133040308
$ws.Range("L12").Value = 1.040282034474445
$ws.Range("M12").Value = 1.059417918388518

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.037894313455375
$ws.Range("D13").Value = 1.047878509338684
$ws.Range("E13").Value = 1.036712559129195
$ws.Range("F13").Value = 1.055924315878616
$ws.Range("I13").Value = 1.045307637544308
$ws.Range("J13").Value = 1.044508530725691
$ws.Range("K13").Value = 1.051437004273414
$ws.Range("L13").Value = 1.040312329530453
$ws.Range("M13").Value = 1.05945366666068

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.038049598576195
$ws.Range("D14").Value = 1.048004761345471
$ws.Range("E14").Value = 1.03684390450409
$ws.Range("F14").Value = 1.05607291275195
$ws.Range("I14").Value = 1.045356321919549
$ws.Range("J14").Value = 1.044603123544608
$ws.Range("K14").Value = 1.051531596446035
$ws.Range("L14").Value = 1.040411599724693
$ws.Range("M14").Value = 1.059570803661023

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.038145288689746
$ws.Range("D15").Value = 1.048082560937926
$ws.Range("E15").Value = 1.036924851862194
$ws.Range("F15").Value = 1.05616448799438
$ws.Range("I15").Value = 1.045386293703875
$ws.Range("J15").Value = 1.044661403197911
$ws.Range("K15").Value = 1.051589872758619
$ws.Range("L15").Value = 1.040472768086928
$ws.Range("M15").Value = 1.059642979536586

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.038702456634549
$ws.Range("D16").Value = 1.048535563571905
$ws.Range("E16").Value = 1.03739632200994
$ws.Range("F16").Value = 1.056697796614182
$ws.Range("I16").Value = 1.045560371387418
$ws.Range("J16").Value = 1.04500058299274
$ws.Range("K16").Value = 1.051928988128898
$ws.Range("L16").Value = 1.040828865532294
$ws.Range("M16").Value = 1.060063132880838

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.039052137387636
$ws.Range("D17").Value = 1.048819873336705
$ws.Range("E17").Value = 1.037692344665592
$ws.Range("F17").Value = 1.057032590881072
$ws.Range("I17").Value = 1.045669238332236
$ws.Range("J17").Value = 1.04521331156393
$ws.Range("K17").Value = 1.052141636616666
$ws.Range("L17").Value = 1.04105229771943
$ws.Range("M17").Value = 1.060326734412797

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.039256166159455
$ws.Range("D18").Value = 1.048985761108059
$ws.Range("E18").Value = 1.037865110743881
$ws.Range("F18").Value = 1.057227965265747
$ws.Range("I18").Value = 1.045732620365364
$ws.Range("J18").Value = 1.045337381386543
$ws.Range("K18").Value = 1.052265645409113
$ws.Range("L18").Value = 1.041182643767281
$ws.Range("M18").Value = 1.060480506456964

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.039325745837501
$ws.Range("D19").Value = 1.049042333780381
$ws.Range("E19").Value = 1.037924036606468
$ws.Range("F19").Value = 1.057294598868913
$ws.Range("I19").Value = 1.045754211928012
$ws.Range("J19").Value = 1.045379684063928
$ws.Range("K19").Value = 1.052307924845999
$ws.Range("L19").Value = 1.041227092067469
$ws.Range("M19").Value = 1.060532941715257

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.039014613126026
$ws.Range("D20").Value = 1.048789363920767
$ws.Range("E20").Value = 1.037660573770256
$ws.Range("F20").Value = 1.056996660832535
$ws.Range("I20").Value = 1.045657570161649
$ws.Range("J20").Value = 1.045190488947107
$ws.Range("K20").Value = 1.052118824072615
$ws.Range("L20").Value = 1.041028323293191
$ws.Range("M20").Value = 1.06029845061913

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.03800386809598
$ws.Range("D21").Value = 1.047967580864273
$ws.Range("E21").Value = 1.036805222158511
$ws.Range("F21").Value = 1.056029150547686
$ws.Range("I21").Value = 1.045341990638374
$ws.Range("J21").Value = 1.044575268753482
$ws.Range("K21").Value = 1.0515037424581
$ws.Range("L21").Value = 1.040382366141418
$ws.Range("M21").Value = 1.059536308918259

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.037368975190798
$ws.Range("D22").Value = 1.047451396879707
$ws.Range("E22").Value = 1.036268350269638
$ws.Range("F22").Value = 1.055421701783377
$ws.Range("I22").Value = 1.045142513176554
$ws.Range("J22").Value = 1.044188362610634
$ws.Range("K22").Value = 1.051116794913299
$ws.Range("L22").Value = 1.039976432309079
$ws.Range("M22").Value = 1.059057289495435

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.037705486380273
$ws.Range("D23").Value = 1.047724987509291
$ws.Range("E23").Value = 1.036552868320497
$ws.Range("F23").Value = 1.055743639471677
$ws.Range("I23").Value = 1.045248360199719
$ws.Range("J23").Value = 1.044393477299317
$ws.Range("K23").Value = 1.051321943710738
$ws.Range("L23").Value = 1.040191605670057
$ws.Range("M23").Value = 1.059311210240929

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.039031568521126
$ws.Range("D24").Value = 1.048803149643752
$ws.Range("E24").Value = 1.037674929358194
$ws.Range("F24").Value = 1.057012895782917
$ws.Range("I24").Value = 1.04566284287102
$ws.Range("J24").Value = 1.045200801539846
$ws.Range("K24").Value = 1.052129132158154
$ws.Range("L24").Value = 1.041039156237118
$ws.Range("M24").Value = 1.060311230793398

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.040572993677315
$ws.Range("D25").Value = 1.050056436651793
$ws.Range("E25").Value = 1.038980977741183
$ws.Range("F25").Value = 1.05848950014274
$ws.Range("I25").Value = 1.046139187070061
$ws.Range("J25").Value = 1.046137214785243
$ws.Range("K25").Value = 1.053064825876424
$ws.Range("L25").Value = 1.042023543558615
$ws.Range("M25").Value = 1.061472389461825
